$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# are first explicitly formatted as Text ("@") so the literal string is kept,
# matching the "Price" column which always stores values as text.
$textCells = @("D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D16", "D18", "D19", "D22", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "25.732.75"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "1.818.31"
$ws.Range("E3").Value = "  -2.71%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "278.88"
$ws.Range("E5").Value = "  -7.27%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "0.5090"
$ws.Range("E7").Value = "  -4.54%  "
$ws.Range("D8").Value = "0.3537"
$ws.Range("E8").Value = "  -5.52%  "
$ws.Range("D9").Value = "44.42"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "0.06691"
$ws.Range("E10").Value = "  -6.88%  "
$ws.Range("D11").Value = "19.94"
$ws.Range("E11").Value = "  -7.90%  "
$ws.Range("D12").Value = "0.8273"
$ws.Range("E12").Value = "  -6.94%  "
$ws.Range("D13").Value = "0.07871"
$ws.Range("E13").Value = "  -3.39%  "
$ws.Range("D14").Value = "1.816.29"
$ws.Range("E14").Value = "  -5.39%  "
$ws.Range("D15").Value = "5.082"
$ws.Range("E15").Value = "  -4.00%  "
$ws.Range("D16").Value = "87.72"
$ws.Range("E16").Value = "  -5.54%  "
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "14.10"
$ws.Range("E18").Value = "  -4.87%  "
$ws.Range("D19").Value = "0.000008044"
$ws.Range("E19").Value = "  -5.36%  "
$ws.Range("D21").Value = "25.785.31"
$ws.Range("E21").Value = "  -3.78%  "
$ws.Range("D22").Value = "4.758"
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("E23").Value = "  -5.80%  "
$ws.Range("D24").Value = "6.109"
$ws.Range("E24").Value = "  -4.41%  "
$ws.Range("D25").Value = "2.244"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("D27").Value = "1.671"
$ws.Range("E27").Value = "  -3.70%  "
$ws.Range("D28").Value = "17.14"
$ws.Range("E28").Value = "  -4.84%  "
$ws.Range("D29").Value = "109.30"
$ws.Range("E29").Value = "  -4.10%  "
$ws.Range("D30").Value = "4.344"
$ws.Range("E30").Value = "  -7.97%  "
$ws.Range("D31").Value = "4.232"
$ws.Range("E31").Value = "  -8.51%  "
$ws.Range("D32").Value = "0.08766"
$ws.Range("E32").Value = "  -4.11%  "
$ws.Range("D33").Value = "0.04892"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").Value = "0.7291"
$ws.Range("E34").Value = "  -9.54%  "
$ws.Range("E35").Value = "  -3.00%  "
$ws.Range("E36").Value = "  -1.58%  "
$ws.Range("D37").Value = "3.159"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("B38").Value = "Frax"
$ws.Range("C38").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D38").Value = "1.000"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "2.361"
$ws.Range("E39").Value = "  -11.94%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "0.01856"
$ws.Range("E40").Value = "  -4.91%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5165"
$ws.Range("E41").Value = "  -15.67%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "0.9707"
$ws.Range("E42").Value = "  -9.20%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "114.37"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "6.224"
$ws.Range("E44").Value = "  -4.50%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "8.016"
$ws.Range("E45").Value = "  -8.69%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.4539"
$ws.Range("E47").Value = "  -12.77%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1370"
$ws.Range("E48").Value = "  -8.22%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "36.48"
$ws.Range("E49").Value = "  -3.14%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.197"
$ws.Range("E50").Value = "  -7.58%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "1.501"
$ws.Range("E51").Value = "  -9.03%  "
